$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 3 and row 4 for columns D, K, L, M, O, P
$cols = @("D", "K", "L", "M", "O", "P")

foreach ($col in $cols) {
    $cell3 = $ws.Range($col + "3")
    $cell4 = $ws.Range($col + "4")
    $temp = $cell3.Value2
    $cell3.Value = $cell4.Value2
    $cell4.Value = $temp
}
